$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The two user-entered inputs on the "Fin Buff Calc" sheet:
#   D3 = Gross Expenditures (From 502 Part C)
#   D5 = Total Labor Cost (From 502 Part L)
# All other figures on the sheet (D6 Labor %, D8 Suggested Standard Added
# Cost, E8 Technician Hours, D9/E9 the 1.5x "detailed" added cost and
# hours) are formulas that depend on D3/D5 and recalculate automatically.
$ws.Range("D3").Value = 100009.74
$ws.Range("D5").Value = 41420.33

$excel.CalculateFullRebuild()
